# Group Import: temporarily allow empty values for column L (Enable Web Hooks).
# CVDLS-220

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# L3 used to hold a literal FALSE; blank it out (will default to NO / FALSE on import)
# and flag it with a light highlight fill so it's easy to spot in the sheet.
$ws.Range("L3").ClearContents()
$ws.Range("L3").Interior.ThemeColor = 10
$ws.Range("L3").Interior.TintAndShade = 0.79998168889431442

# L11 had been left blank/mis-filled; put back the marker text so QA remembers to
# reset it to empty again later.
$ws.Range("L11").Value = "CHANGE BACK TO EMPTY"

# Give column K (11) some breathing room now that it's visible in the selection flow.
$ws.Columns.Item(11).ColumnWidth = 21.830729166666668

# Leave the cursor parked on the cell we just cleared.
[void]$ws.Range("L3").Select()
